$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue "D2" '58.224.48'
Set-TextValue "E2" '  -0.37%  '
Set-TextValue "D3" '2.582.10'
Set-TextValue "E3" '  -2.38%  '
Set-TextValue "D4" '1.00'
Set-TextValue "E4" '  -0.01%  '
Set-TextValue "D5" '518.17'
Set-TextValue "E5" '  -0.56%  '
Set-TextValue "D6" '142.87'
Set-TextValue "E6" '  -0.83%  '
Set-TextValue "E7" '  -0.12%  '
Set-TextValue "E8" '  -0.40%  '
Set-TextValue "D9" '2.599.32'
Set-TextValue "E9" '  -2.04%  '
Set-TextValue "D10" '6.85'
Set-TextValue "E10" '  +2.42%  '
Set-TextValue "E11" '  -1.49%  '
Set-TextValue "E12" '  -3.50%  '
Set-TextValue "E13" '  -1.15%  '
Set-TextValue "D14" '3.041.23'
Set-TextValue "E14" '  -2.17%  '
Set-TextValue "D15" '58.136.01'
Set-TextValue "E15" '  -0.51%  '
Set-TextValue "D16" '20.33'
Set-TextValue "E16" '  -2.34%  '
Set-TextValue "E17" '  -1.71%  '
Set-TextValue "D18" '2.593.43'
Set-TextValue "E18" '  -1.91%  '
Set-TextValue "D19" '344.02'
Set-TextValue "E19" '  +1.98%  '
Set-TextValue "E20" '  -2.13%  '
Set-TextValue "D21" '10.28'
Set-TextValue "E21" '  -1.75%  '
Set-TextValue "D22" '6.34'
Set-TextValue "E22" '  +0.72%  '
Set-TextValue "E23" '  +0.13%  '
Set-TextValue "D24" '66.28'
Set-TextValue "E24" '  +2.83%  '
Set-TextValue "E25" '  -1.22%  '
Set-TextValue "E26" '  -5.19%  '
Set-TextValue "E27" '  -0.24%  '
Set-TextValue "D28" '2.697.05'
Set-TextValue "E28" '  -2.41%  '
Set-TextValue "E29" '  -1.42%  '
Set-TextValue "E30" '  -6.53%  '
Set-TextValue "E31" '  -0.01%  '
Set-TextValue "D32" '6.22'
Set-TextValue "E32" '  -6.03%  '
Set-TextValue "E33" '  -0.50%  '
Set-TextValue "D34" '18.76'
Set-TextValue "E34" '  -0.39%  '
Set-TextValue "D35" '149.65'
Set-TextValue "E35" '  -2.10%  '
Set-TextValue "E36" '  -2.54%  '
Set-TextValue "E37" '  -2.98%  '
Set-TextValue "D38" '0.873'
Set-TextValue "E38" '  -3.75%  '
Set-TextValue "D39" '0.837'
Set-TextValue "E39" '  -2.14%  '
Set-TextValue "D40" '35.93'
Set-TextValue "E40" '  -2.37%  '
Set-TextValue "D41" '1.45'
Set-TextValue "E41" '  +0.33%  '
Set-TextValue "D42" '3.54'
Set-TextValue "E42" '  -2.49%  '
Set-TextValue "E43" '  -0.20%  '
Set-TextValue "D44" '273.88'
Set-TextValue "E44" '  +1.75%  '
Set-TextValue "E45" '  -2.77%  '
Set-TextValue "D46" '10.64'
Set-TextValue "E46" '  +0.02%  '
Set-TextValue "D47" '0.0952'
Set-TextValue "E47" '  -1.73%  '
Set-TextValue "B48" 'EnergySwap'
Set-TextValue "C48" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D48" '18.87'
Set-TextValue "E48" '  -2.65%  '
Set-TextValue "B49" 'Hedera'
Set-TextValue "C49" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D49" '0.0525'
Set-TextValue "E49" '  -2.03%  '
Set-TextValue "D50" '1.973.09'
Set-TextValue "E50" '  -3.45%  '
Set-TextValue "B51" 'InjectiveProtocol'
Set-TextValue "C51" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D51" '18.58'
Set-TextValue "E51" '  +1.45%  '
